# edit.ps1 - "First pass at ch3 DCAF Components"
#
#  1. Update the cached "datetimeFigureOut" date placeholder text from
#     6/6/2018 -> 2018-06-18 on the slide master, every slide layout and
#     the notes master.
#  2. Insert three new content slides (positions 3, 4, 5) that were moved
#     in while putting together the new "DCAF Components" chapter-3
#     session. The middle slide ("Which Template to Choose") is marked
#     hidden.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Date placeholder text fix-up (6/6/2018 -> 2018-06-18)
# ---------------------------------------------------------------------
$newDate = "2018-06-18"

$p.SlideMaster.Shapes.Item("Date Placeholder 3").TextFrame.TextRange.Text = $newDate

for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$p.NotesMaster.Shapes.Item("Date Placeholder 2").TextFrame.TextRange.Text = $newDate

# ---------------------------------------------------------------------
# 2. New slide 3: "What type of module do I need"
# ---------------------------------------------------------------------
$s3 = $p.Slides.Add(3, 2)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = " What type of module do I need"

$body3 = $s3.Shapes.Item(2).TextFrame.TextRange
$body3.Text = "placeholder`rDo we really need a custom module??`rUse Existing modules`rModify "
$body3.InsertAfter("Existing modules")
$body3.InsertAfter("`rWhen ")
$body3.InsertAfter("do you define the inputs and outputs")
$body3.InsertAfter("`rWhat is the module doing")
$body3.InsertAfter("`rplaceholder")

$body3.Paragraphs(1).Text = ""
$body3.Paragraphs(7).Text = ""
$body3.Paragraphs(3).IndentLevel = 2
$body3.Paragraphs(4).IndentLevel = 2

# ---------------------------------------------------------------------
# 3. New slide 4: "Which Template to Choose" (hidden)
# ---------------------------------------------------------------------
$s4 = $p.Slides.Add(4, 2)
$s4.SlideShowTransition.Hidden = 1
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Which Template to Choose"

$body4 = $s4.Shapes.Item(2).TextFrame.TextRange
$body4.Text = "Dynamic`rUnknown number of Channels.`rFinal user defines the inputs`rModbus`rNeeds more codding`rStatic Template`rNumber of channels decided by module developer`rMostly Scripted`rPID`rplaceholder"

$body4.Paragraphs(10).Text = ""
$body4.Paragraphs(2).IndentLevel = 2
$body4.Paragraphs(3).IndentLevel = 2
$body4.Paragraphs(4).IndentLevel = 2
$body4.Paragraphs(5).IndentLevel = 2
$body4.Paragraphs(7).IndentLevel = 2
$body4.Paragraphs(8).IndentLevel = 2
$body4.Paragraphs(9).IndentLevel = 2

# ---------------------------------------------------------------------
# 4. New slide 5: "When do I define the inputs"
# ---------------------------------------------------------------------
$s5 = $p.Slides.Add(5, 2)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "When do I define the inputs"

$body5 = $s5.Shapes.Item(2).TextFrame.TextRange
$body5.Text = "Configuration `rUnknown number of inputs. Final user defines the inputs`rDynamic Template`rModbus`rplaceholder`rDevelopment`rAlways the same number of inputs`rStatic Template`rPID`rplaceholder"

$body5.Paragraphs(5).Text = ""
$body5.Paragraphs(10).Text = ""
$body5.Paragraphs(2).IndentLevel = 2
$body5.Paragraphs(3).IndentLevel = 2
$body5.Paragraphs(4).IndentLevel = 2
$body5.Paragraphs(7).IndentLevel = 2
$body5.Paragraphs(8).IndentLevel = 2
$body5.Paragraphs(9).IndentLevel = 2
